# The commit drops the first event row (上饶·原×铁×崩only（取消）) from both
# the "展览" and "全部类型" sheets. The leading "序号" column (A) keeps its
# original 0-based numbering per physical row, but every event's B:I data
# (date/name/place/time/count/price/link/cover) shifts up by one row, and
# the final row (32, 南昌·代号鸢盛花行only) disappears, shrinking the used
# range from A1:I32 to A1:I31. A handful of "想去人数" (F column) interest
# counts were also bumped for some of the remaining events (re-scraped at a
# later time), so those are patched in afterwards.

$wb = $excel.ActiveWorkbook

# F-column ("想去人数") values to apply after the shift, keyed by row number.
$fUpdates = @{
    2  = 3136
    6  = 1743
    8  = 97
    11 = 1437
    13 = 562
    14 = 357
    15 = 76
    16 = 13
    17 = 80
    18 = 66
    21 = 96
    23 = 3372
    24 = 407
    25 = 290
    26 = 465
    27 = 46
    30 = 1114
    31 = 116
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Shift every event's B:I data up by one row (rows 3..32 -> 2..31),
    # leaving column A (the fixed 0-based index) untouched.
    $src = $ws.Range("B3:I32")
    $dst = $ws.Range("B2")
    $src.Copy($dst)

    # Drop the now-duplicated last row so the used range becomes A1:I31.
    $ws.Rows.Item(32).Delete()

    foreach ($rowNum in $fUpdates.Keys) {
        $ws.Range("F$rowNum").Value = $fUpdates[$rowNum]
    }
}
